# Update "想去人数" (wanted-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 5311
$ws1.Range("F12").Value = 146
$ws1.Range("F15").Value = 370
$ws1.Range("F18").Value = 182
$ws1.Range("F21").Value = 6089
$ws1.Range("F22").Value = 6089
$ws1.Range("F26").Value = 7037
$ws1.Range("F29").Value = 3265
$ws1.Range("F34").Value = 140
$ws1.Range("F40").Value = 933
$ws1.Range("F41").Value = 1157

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 5311
$ws4.Range("F15").Value = 146
$ws4.Range("F18").Value = 370
$ws4.Range("F22").Value = 182
$ws4.Range("F25").Value = 6089
$ws4.Range("F29").Value = 7037
$ws4.Range("F32").Value = 3265
$ws4.Range("F38").Value = 140
$ws4.Range("F44").Value = 933
$ws4.Range("F45").Value = 1157
